$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "TestValue"
Write-Host "Sheet name:" $ws.Name
Write-Host "Sheet count:" $wb.Worksheets.Count
